$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the answer value for row 2 first (so "a" gets the earlier shared-string index)
$ws.Range("M2").Value = "a"

# Add new "Correct Option" column header in M1, matching the style of the other headers
$ws.Range("M1").Value = "Correct Option"
$ws.Range("M1").Interior.ColorIndex = 6

# Resize column M to fit its content (matches the bestFit width of the other columns)
$ws.Columns.Item(13).ColumnWidth = 13

# Update the active selection to N1, matching the post-edit state
$ws.Range("N1").Select()
